$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename shared strings used in column B (filename) - strip "_CRR_DiadFit" suffix
$ws.Range("B4").Value = "K21-64-FI#1"
$ws.Range("B5").Value = "K21-64-FI#3"
$ws.Range("B7").Value = "K21-64-FI#5"
$ws.Range("B8").Value = "K21-64-FI#6"
$ws.Range("B9").Value = "K21-64-FI#7"

# Update numeric cell values
$ws.Range("C2").Value = 103.3001276631323
$ws.Range("D2").Value = 0.003265494907594278
$ws.Range("E2").Value = 0.002477130814084896
$ws.Range("F2").Value = 1285.470024067249
$ws.Range("G2").Value = 0.002266959001439887
$ws.Range("H2").Value = 988.9448337599146
$ws.Range("I2").Value = 1285.470074069749
$ws.Range("J2").Value = 1433.106807154252
$ws.Range("K2").Value = 0.552660665644079
$ws.Range("L2").Value = 2.918337680969268
$ws.Range("M2").Value = 0.5834875377734332
$ws.Range("N2").Value = 1.105321331288158
$ws.Range("Q2").Value = 0.0009985359061543907
$ws.Range("V2").Value = 0
$ws.Range("C3").Value = 103.3338104977101
$ws.Range("D3").Value = 0.003339245779126
$ws.Range("E3").Value = 0.002441025094158034
$ws.Range("F3").Value = 1285.409241475892
$ws.Range("G3").Value = 0.002107414120404731
$ws.Range("H3").Value = 815.58347963966
$ws.Range("I3").Value = 1285.409291478392
$ws.Range("J3").Value = 1225.330007618777
$ws.Range("K3").Value = 0.5569708025234513
$ws.Range("L3").Value = 2.339882415644469
$ws.Range("M3").Value = 0.6539937961029044
$ws.Range("N3").Value = 1.113941605046903
$ws.Range("Q3").Value = 0.001231831658721269
$ws.Range("V3").Value = 0
$ws.Range("D4").Value = 0.004350429593681681
$ws.Range("E4").Value = 0.003332763968836868
$ws.Range("G4").Value = 0.003081911230931819
$ws.Range("Q4").Value = 0.001268518362749862
$ws.Range("V4").Value = 0
$ws.Range("C5").Value = 103.3010339772266
$ws.Range("D5").Value = 0.002899358026653716
$ws.Range("E5").Value = 0.002119282776369108
$ws.Range("F5").Value = 1285.380956189148
$ws.Range("G5").Value = 0.001829297903040009
$ws.Range("H5").Value = 1350.101643341807
$ws.Range("I5").Value = 1285.381006191648
$ws.Range("J5").Value = 2155.06232825577
$ws.Range("K5").Value = 0.586194237919843
$ws.Range("L5").Value = 2.606146132236945
$ws.Range("M5").Value = 0.6770808892749949
$ws.Range("N5").Value = 1.172388475839686
$ws.Range("Q5").Value = 0.001070060123613707
$ws.Range("V5").Value = 0
$ws.Range("AA5").Value = 1265.141738155238
$ws.Range("AB5").Value = 225.1444621488707
$ws.Range("AC5").Value = 0.7535648176856447
$ws.Range("D6").Value = 0.002717581954345735
$ws.Range("E6").Value = 0.002120968633252415
$ws.Range("G6").Value = 0.001993567356675061
$ws.Range("Q6").Value = 0.0007240145976706735
$ws.Range("V6").Value = 0
$ws.Range("C7").Value = 103.3060250519686
$ws.Range("D7").Value = 0.0040037347721505
$ws.Range("E7").Value = 0.003100328033770072
$ws.Range("F7").Value = 1285.37382745755
$ws.Range("G7").Value = 0.002895480236843717
$ws.Range("H7").Value = 999.1651985774955
$ws.Range("I7").Value = 1285.37382745755
$ws.Range("J7").Value = 1571.528012032771
$ws.Range("K7").Value = 0.5895992431694473
$ws.Range("L7").Value = 4.098160685460347
$ws.Range("M7").Value = 0.6258550187037439
$ws.Range("N7").Value = 1.179198486338895
$ws.Range("Q7").Value = 0.001108254535306783
$ws.Range("V7").Value = 0
$ws.Range("D8").Value = 0.003400856191057229
$ws.Range("E8").Value = 0.002592957819105195
$ws.Range("G8").Value = 0.002386182729319969
$ws.Range("H8").Value = 741.9066878607191
$ws.Range("K8").Value = 0.5933826513295981
$ws.Range("L8").Value = 2.548626843598762
$ws.Range("Q8").Value = 0.001014673461737261
$ws.Range("V8").Value = 0
$ws.Range("C9").Value = 103.3077444099038
$ws.Range("D9").Value = 0.003283041856616572
$ws.Range("E9").Value = 0.002432674674102685
$ws.Range("F9").Value = 1285.341648299921
$ws.Range("G9").Value = 0.002155682456367518
$ws.Range("H9").Value = 1184.564148378073
$ws.Range("I9").Value = 1285.341698302421
$ws.Range("J9").Value = 1945.250843440412
$ws.Range("K9").Value = 0.6050003855460366
$ws.Range("L9").Value = 2.982833638240457
$ws.Range("M9").Value = 0.6693677755212957
$ws.Range("N9").Value = 1.210000771092073
$ws.Range("Q9").Value = 0.001127359400249055
$ws.Range("V9").Value = 0
$ws.Range("AA9").Value = 1265.101042386816
$ws.Range("AB9").Value = 219.5095786259168
$ws.Range("AC9").Value = 0.8616414706285049

# Clear cells that should no longer have values
$ws.Range("AA7").ClearContents()
$ws.Range("AB7").ClearContents()
$ws.Range("AC7").ClearContents()
